# Update countries & provincias Spain
#
# The source feed was re-pulled (20 Mar 2020, 06:16 instead of 05:46). Several
# countries' case counts grew, which re-sorted the (descending-by-total-cases)
# table: a handful of rows keep their row number but now show a different
# country with refreshed figures, while the rest of that country's old figures
# slide down to the next row. Only the rows whose effective country/number
# content actually changes are touched below; everything else is left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp
$ws.Range("A1").Value = "Datos actualizados a 20 de Marzo de 2020 a las 06:16"

$rows = @(
    @{Row=9;   Name="Estados Unidos";      B=14354; C=565; D=125; E=14012; F=64; G=9; H=217},
    @{Row=23;  Name="Australia";           B=814;   C=58;  D=46;  E=761;   F=1;  G=0; H=7},
    @{Row=24;  Name="Portugal";            B=786;   C=0;   D=4;   E=778;   F=20; G=0; H=4},
    @{Row=40;  Name="Tailandia";           B=322;   C=50;  D=42;  E=279;   F=1;  G=0; H=1},
    @{Row=41;  Name="Eslovenia";           B=319;   C=0;   D=0;   E=318;   F=6;  G=0; H=1},
    @{Row=42;  Name="Indonesia";           B=309;   C=0;   D=15;  E=269;   F=0;  G=0; H=25},
    @{Row=43;  Name="Barein";              B=279;   C=0;   D=110; E=168;   F=4;  G=0; H=1},
    @{Row=44;  Name="Rumania";             B=277;   C=0;   D=25;  E=252;   F=5;  G=0; H=0},
    @{Row=45;  Name="Arabia Saudita";      B=274;   C=0;   D=8;   E=266;   F=0;  G=0; H=0},
    @{Row=53;  Name="India";               B=197;   C=3;   D=20;  E=173;   F=0;  G=0; H=4},
    @{Row=72;  Name="Costa Rica";          B=89;    C=2;   D=0;   E=87;    F=2;  G=1; H=2},
    @{Row=124; Name="Polinesia Francesa";  B=11;    C=5;   D=0;   E=11;    F=0;  G=0; H=0},
    @{Row=125; Name="Ruanda";              B=11;    C=0;   D=0;   E=11;    F=0;  G=0; H=0},
    @{Row=126; Name="Monaco";              B=10;    C=0;   D=0;   E=10;    F=0;  G=0; H=0},
    @{Row=127; Name="Gibraltar";           B=10;    C=0;   D=2;   E=8;     F=0;  G=0; H=0},
    @{Row=128; Name="Trinidad yTobago";    B=9;     C=0;   D=0;   E=9;     F=0;  G=0; H=0},
    @{Row=129; Name="Etiopia";             B=9;     C=2;   D=0;   E=9;     F=0;  G=0; H=0},
    @{Row=130; Name="Guatemala";           B=9;     C=0;   D=0;   E=8;     F=0;  G=0; H=1},
    @{Row=131; Name="Costa de Marfil";     B=9;     C=0;   D=1;   E=8;     F=0;  G=0; H=0},
    @{Row=132; Name="Mauricio";            B=7;     C=0;   D=0;   E=7;     F=0;  G=0; H=0},
    @{Row=133; Name="Kenia";               B=7;     C=0;   D=0;   E=7;     F=0;  G=0; H=0},
    @{Row=134; Name="Puerto Rico";         B=6;     C=0;   D=0;   E=6;     F=0;  G=0; H=0},
    @{Row=135; Name="Kirguistan";          B=6;     C=3;   D=0;   E=6;     F=0;  G=0; H=0},
    @{Row=136; Name="Seychelles";          B=6;     C=0;   D=0;   E=6;     F=0;  G=0; H=0},
    @{Row=137; Name="Guinea Ecuatorial";   B=6;     C=0;   D=0;   E=6;     F=0;  G=0; H=0},
    @{Row=138; Name="Mongolia";            B=6;     C=0;   D=0;   E=6;     F=0;  G=0; H=0},
    @{Row=139; Name="Tanzania";            B=6;     C=0;   D=0;   E=6;     F=0;  G=0; H=0}
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Name
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
    $ws.Cells.Item($r.Row, 8).Value = $r.H
}
